$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 213, pushing existing rows 213:313 down to 214:314.
$ws.Rows.Item(213).Insert()

# Populate the newly inserted row 213 with the new data record.
$ws.Range("A213").Value = 7
$ws.Range("B213").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C213").Value = "Ñuble"
$ws.Range("D213").Value = 45016
$ws.Range("E213").Value = 16
$ws.Range("F213").Value = 100112017
$ws.Range("G213").Value = "Apio"
$ws.Range("H213").Value = "Americana (o)"
$ws.Range("I213").Value = "Segunda"
$ws.Range("J213").Value = 120
$ws.Range("K213").Value = 6000
$ws.Range("L213").Value = 6000
$ws.Range("M213").Value = 6000
$ws.Range("N213").Value = '$/docena de matas'
$ws.Range("O213").Value = "Provincia del Elquí"
$ws.Range("P213").Value = 1000
$ws.Range("Q213").Value = 6
$ws.Range("R213").Value = "Hortaliza"
